$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Tom" dataset group (rows 74-85) gets consolidated into the new
# "Longitudinal analysis ..." (PMP22 / CMT1A) dataset group (rows 74-81).
# ---------------------------------------------------------------------------

# First, copy cell formatting onto the rows that will host the new data so
# every row ends up visually consistent with its sibling rows:
#  - rows 76-77 should look like rows 74-75 (tan "title" style)
#  - rows 78-81 should look like the old rows 80-81 ("WT"/no-fill style)
$ws.Range("A74:D74").Copy()
$ws.Range("A76:D77").PasteSpecial(-4122)

$ws.Range("A80:D80").Copy()
$ws.Range("A78:D81").PasteSpecial(-4122)

# Now write the new values (8 rows instead of the original 12).
$titleTg = "Longitudinal analysis of lipid changes in the sciatic nerve caused by overexpression of PMP22 in models of CMT1A"
$titleNoTg = $titleTg + " (no TG)"

$ws.Range("A74").Value = "NLA_913"
$ws.Range("B74").Value = $titleTg
$ws.Range("C74").Value = "C22"
$ws.Range("D74").Value = "NLA_091"

$ws.Range("A75").Value = "NLA_913"
$ws.Range("B75").Value = $titleTg
$ws.Range("C75").Value = "C22 – WT"
$ws.Range("D75").Value = "NLA_092"

$ws.Range("A76").Value = "NLA_913"
$ws.Range("B76").Value = $titleTg
$ws.Range("C76").Value = "C3"
$ws.Range("D76").Value = "NLA_093"

$ws.Range("A77").Value = "NLA_913"
$ws.Range("B77").Value = $titleTg
$ws.Range("C77").Value = "C3 – WT"
$ws.Range("D77").Value = "NLA_094"

$ws.Range("A78").Value = "NLA_914"
$ws.Range("B78").Value = $titleNoTg
$ws.Range("C78").Value = "C22 (no TG)"
$ws.Range("D78").Value = "NLA_095"

$ws.Range("A79").Value = "NLA_914"
$ws.Range("B79").Value = $titleNoTg
$ws.Range("C79").Value = "C22 – WT (no TG)"
$ws.Range("D79").Value = "NLA_096"

$ws.Range("A80").Value = "NLA_914"
$ws.Range("B80").Value = $titleNoTg
$ws.Range("C80").Value = "C3 (no TG)"
$ws.Range("D80").Value = "NLA_097"

$ws.Range("A81").Value = "NLA_914"
$ws.Range("B81").Value = $titleNoTg
$ws.Range("C81").Value = "C3 – WT (no TG)"
$ws.Range("D81").Value = "NLA_098"

# The remaining 4 rows of the old layout are no longer needed.
$ws.Rows("82:85").Delete()

# ---------------------------------------------------------------------------
# Cosmetic / view updates that came along with this edit.
# ---------------------------------------------------------------------------
$ws.Columns(2).ColumnWidth = 105.46

$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("B81").Select()
